# FBL-21 - Import excel data & create player, staff templates
# Turns the single "staff roster" sheet into a two-sheet template:
#   Sheet1 = instructions / legend
#   Sheet2 = data-entry sheet (No / Name / Role / Nationality) with a
#            decimal data-validation rule on the Role column.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "Sheet1"

# --- Sheet1: wipe the old single-row header, it gets replaced wholesale ---
$ws1.Range("A1:F1").Clear()

# Row 1-2: bold intro text (column A only)
$ws1.Range("A1").Value = "Đây là trang hướng dẫn"
$ws1.Range("A1").Font.Bold = $true
$ws1.Range("A2").Value = "Sheet2 là trang chứa dữ liệu, khi submit, không được đổi tên Sheet2"
$ws1.Range("A2").Font.Bold = $true

# Row 4: plain (non-bold) section title
$ws1.Range("A4").Value = "Ý nghĩa các trường thông tin như sau"

# Rows 5-8: legend table, column A = field name (bold + border), column B = meaning (border)
$ws1.Range("A5").Value = "No"
$ws1.Range("B5").Value = "Số thứ tự"

$ws1.Range("A6").Value = "Name"
$ws1.Range("B6").Value = "Tên thành viên"

$ws1.Range("A7").Value = "Role"
$ws1.Range("B7").Value = "Loại thành viên"

$ws1.Range("A8").Value = "Nationality"
$ws1.Range("B8").Value = "Quốc tịch của cầu thủ"

$ws1.Range("A5:A8").Font.Bold = $true
$ws1.Range("A5:B8").Borders.LineStyle = 1

# Column widths for the new layout
$ws1.Columns.Item(1).ColumnWidth = 13.14
$ws1.Columns.Item(2).ColumnWidth = 55.57

# --- Sheet2: brand new data sheet, placed right after Sheet1 ---
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Sheet2"

$ws2.Range("A1").Value = "No"
$ws2.Range("B1").Value = "Name"
$ws2.Range("C1").Value = "Role"
$ws2.Range("D1").Value = "Nationality"
$ws2.Range("A1:D1").Font.Bold = $true
$ws2.Range("F1").Font.Bold = $true

$ws2.Range("A2").Value = 1
$ws2.Range("B2").Value = "Pep"
$ws2.Range("C2").Value = 0
$ws2.Range("D2").Value = "Tây Ban Nha"

$ws2.Range("A3").Value = 2
$ws2.Range("B3").Value = "Trợ lí Pep"
$ws2.Range("C3").Value = 1
$ws2.Range("D3").Value = "Tây Ban Nha"

$ws2.Range("C1").NumberFormat = "0"
$ws2.Range("C2:C3").NumberFormat = "0"

# Column widths
$ws2.Columns.Item(2).ColumnWidth = 20.57
$ws2.Columns.Item(3).ColumnWidth = 17.71
$ws2.Columns.Item(4).ColumnWidth = 22
$ws2.Columns.Item(6).ColumnWidth = 16.86

# Data validation: Role must be a decimal 0-2 (0 = head coach, 1 = assistant, 2 = other)
$validation = $ws2.Range("C2:C1048576").Validation
$validation.Add(2, 1, 1, 0, 2)
$validation.IgnoreBlank = $true
$validation.ShowInput = $true
$validation.ShowError = $true
$validation.InputTitle = "Hướng dẫn"
$validation.InputMessage = "0 là HLV trưởng,1 là trợ lí HLV, 2 là nhân viên khác, nếu có nhiều hơn thành viên là HLV trưởng, mặc định người đầu tiên sẽ được chọn"

[void]$ws2.Range("B4").Select()

# Sheet1 is the tab that ends up selected/active in the saved workbook
[void]$ws1.Activate()
[void]$ws1.Range("A9").Select()

Write-Output "done"
